$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that currently sits
#    right after the H1 title paragraph.
# ------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Meta description:*") {
        $p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 2. Find the paragraph that still holds the old DALLE image prompt
#    (the very last paragraph of the document).
# ------------------------------------------------------------------
$dalleIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Prompt for DALLE:*") {
        $dalleIndex = $i
        break
    }
}

$dallePara = $d.Paragraphs.Item($dalleIndex)

# ------------------------------------------------------------------
# 3. Insert a brand-new paragraph right before it containing the bold
#    title text "Play Deepsea Riches Free - Slot Game Review 2021".
# ------------------------------------------------------------------
$prevPara = $d.Paragraphs.Item($dalleIndex - 1)
$insPoint = $prevPara.Range
$insPoint.Collapse(0)                 # wdCollapseEnd
$insPoint.InsertAfter("Play Deepsea Riches Free - Slot Game Review 2021" + [char]13)

$newTitlePara = $d.Paragraphs.Item($dalleIndex)
$newTitleRange = $newTitlePara.Range
$newTitleRange.MoveEnd(1, -1) | Out-Null   # wdCharacter, exclude paragraph mark
$newTitleRange.Font.Bold = 1

# ------------------------------------------------------------------
# 4. Replace the text of the (now shifted) DALLE paragraph with the
#    meta-description sentence, keeping its italic formatting intact.
# ------------------------------------------------------------------
$dallePara = $d.Paragraphs.Item($dalleIndex + 1)
$dalleRange = $dallePara.Range
$dalleRange.MoveEnd(1, -1) | Out-Null      # wdCharacter, exclude paragraph mark
$dalleRange.Text = "Explore sunken pirate treasures in Deepsea Riches - read our review and play the game for free on your PC or mobile device."
